# Generate Report for Handoff
# - Flip status from "In Translation" to "Ready for handoff" for the zh-cn / de-de
#   rows (Overview sheet status columns + the per-locale "Status" column on each
#   locale sheet).
# - Bump the corresponding handoff timestamps.
# - Widen the status columns so the longer "Ready for handoff" text fits.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-08-30 22:42:59"

$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332

# --- zh-cn sheet ------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-30 22:42:55"

$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333332

# --- de-de sheet ------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-30 22:42:59"

$dede.Columns.Item(3).ColumnWidth = 16.333333333333332
